$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(20250513, 20250513, 20250513, 20250513, 20250513, 20250513, 20250513)
$timepoints = @(0, 1, 2, 3, 4, 5, 6)
$temps = @("18C", "18C", "18C", "42C", "42C", "42C", "42C")

for ($i = 0; $i -lt 7; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $timepoints[$i]
    $ws.Cells.Item($row, 3).Value = $temps[$i]
}

$ws.Range("D15").Select()
